$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 72: correct the date/time value in column A ---
$ws.Cells.Item(72, 1).Value = 45447.2916666667

# --- Row 73: new row appended (date, volume, high, low, open, close, adj_close, ticker) ---
$ws.Cells.Item(73, 1).Value = 45448.3091782407
# Reuse A72's exact formatting (date number format + font) for the new date cell
$ws.Cells.Item(72, 1).Copy()
$ws.Cells.Item(73, 1).PasteSpecial(-4122)

$ws.Cells.Item(73, 2).Value = 1500
$ws.Cells.Item(73, 3).Value = 2.97000002861023
$ws.Cells.Item(73, 4).Value = 2.97000002861023
$ws.Cells.Item(73, 5).Value = 2.97000002861023
$ws.Cells.Item(73, 6).Value = 2.97000002861023

# Column G (adj_close) is stored as text in this workbook, so force the
# numeric-looking string to stay text with a leading apostrophe, then drop
# back to the default "Normal" style (no special text number format).
$ws.Cells.Item(73, 7).Value = "'2.97000002861023"
$ws.Cells.Item(73, 7).Style = "Normal"

$ws.Cells.Item(73, 8).Value = "ESPE.MI"
